{"js": "// Update the date line and the 25 division-problem answer cells in the\n// single table, preserving each run/paragraph's existing formatting by\n// replacing text through the existing paragraph/cell range rather than\n// rebuilding runs from scratch.\n\n// 1) Update the date paragraph (the first paragraph in the document body,\n//    before the table).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange().insertText(\"2023-09-14 Thursday\", Word.InsertLocation.replace);\n\n// 2) Update the table cells. The table has 20 rows x 5 columns, but only\n//    every 4th row (0, 4, 8, 12, 16) actually holds answer text - the\n//    rows in between are blank spacer rows. Values below are listed in\n//    reading order (row-major) for just those five populated rows.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newRowValues = {\n  0: [\"59\u00f74=14, 3\", \"24\u00f73=8, 0\", \"10\u00f74=2, 2\", \"42\u00f77=6, 0\", \"56\u00f78=7, 0\"],\n  4: [\"71\u00f76=11, 5\", \"98\u00f72=49, 0\", \"91\u00f79=10, 1\", \"15\u00f72=7, 1\", \"75\u00f79=8, 3\"],\n  8: [\"30\u00f73=10, 0\", \"11\u00f79=1, 2\", \"13\u00f77=1, 6\", \"47\u00f73=15, 2\", \"96\u00f75=19, 1\"],\n  12: [\"84\u00f73=28, 0\", \"43\u00f75=8, 3\", \"69\u00f74=17, 1\", \"67\u00f75=13, 2\", \"65\u00f78=8, 1\"],\n  16: [\"49\u00f75=9, 4\", \"16\u00f72=8, 0\", \"92\u00f78=11, 4\", \"40\u00f72=20, 0\", \"20\u00f76=3, 2\"],\n};\n\nfor (const rowIndex of Object.keys(newRowValues)) {\n  const rowIdx = parseInt(rowIndex, 10);\n  const values = newRowValues[rowIndex];\n  for (let colIdx = 0; colIdx < values.length; colIdx++) {\n    const cell = table.getCell(rowIdx, colIdx);\n    cell.body.getRange().insertText(values[colIdx], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-problem answer cells in the\n# single table, preserving each paragraph/run's existing formatting by\n# assigning to the existing Range.Text rather than rebuilding the runs.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph (first paragraph in the document, before\n#    the table).\n$d.Paragraphs.Item(1).Range.Text = \"2023-09-14 Thursday\"\n\n# 2) Update the table cells. The table has 20 rows x 5 columns, but only\n#    every 4th row (1, 5, 9, 13, 17 in 1-based COM indexing) actually\n#    holds answer text - the rows in between are blank spacer rows.\n$t = $d.Tables.Item(1)\n\n$newRowValues = @{\n    1  = @(\"59\u00f74=14, 3\", \"24\u00f73=8, 0\", \"10\u00f74=2, 2\", \"42\u00f77=6, 0\", \"56\u00f78=7, 0\")\n    5  = @(\"71\u00f76=11, 5\", \"98\u00f72=49, 0\", \"91\u00f79=10, 1\", \"15\u00f72=7, 1\", \"75\u00f79=8, 3\")\n    9  = @(\"30\u00f73=10, 0\", \"11\u00f79=1, 2\", \"13\u00f77=1, 6\", \"47\u00f73=15, 2\", \"96\u00f75=19, 1\")\n    13 = @(\"84\u00f73=28, 0\", \"43\u00f75=8, 3\", \"69\u00f74=17, 1\", \"67\u00f75=13, 2\", \"65\u00f78=8, 1\")\n    17 = @(\"49\u00f75=9, 4\", \"16\u00f72=8, 0\", \"92\u00f78=11, 4\", \"40\u00f72=20, 0\", \"20\u00f76=3, 2\")\n}\n\nforeach ($rowIndex in $newRowValues.Keys) {\n    $values = $newRowValues[$rowIndex]\n    for ($colIndex = 1; $colIndex -le $values.Count; $colIndex++) {\n        $cell = $t.Cell($rowIndex, $colIndex)\n        $cell.Range.Text = $values[$colIndex - 1]\n    }\n}\n"}
